$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.9969661396677152
$ws.Range("C3").Value = 0.996946930279147
$ws.Range("D3").Value = 0.9971469744954022

$ws.Range("B4").Value = 0.9976825199215265
$ws.Range("C4").Value = 0.9976810172012281
$ws.Range("D4").Value = 0.997695021204274

$ws.Range("B5").Value = 0.9860838532599849
$ws.Range("C5").Value = 0.9874104151767233
$ws.Range("D5").Value = 0.987983614392475
